$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start clean: wipe all existing cell content/formatting ---
$ws.Cells.Clear()

# --- Write cell values (row/col layout per target) ---
# Row 2 - thin spacer row (no text)
# Row 3 - merged title "Uncertainty scenarios"
$ws.Cells.Item(3,2).Value = "Uncertainty scenarios"

# Row 4-7 - uncertainty scenario definitions
$ws.Cells.Item(4,2).Value = "Base"
$ws.Cells.Item(4,3).Value = "Standard bootstrap-simulation (ommiting steps 5 and 6 that include ageing error and growth variability in the Bootstrap-Simulation framework)"

$ws.Cells.Item(5,2).Value = "AE"
$ws.Cells.Item(5,3).Value = "Bootstrap-simulation including ageing error only"

$ws.Cells.Item(6,2).Value = "GV"
$ws.Cells.Item(6,3).Value = "Bootstrap-simulation including growth variabliity only"

$ws.Cells.Item(7,2).Value = "AE & GV"
$ws.Cells.Item(7,3).Value = "Bootstrap-simulation including both ageing error and growth variability"

# Row 8 - merged "Treatments" sub-header
$ws.Cells.Item(8,2).Value = "Treatments"

# Row 9-11 - treatment definitions
$ws.Cells.Item(9,2).Value = "Growth variaiblity treatment"
$ws.Cells.Item(9,3).Value = "Resample lengths for a given age after pooling age-length data across survey years ('Pooled') or using annual age-length data ('Annual')"

$ws.Cells.Item(10,2).Value = "Length bin treatment"
$ws.Cells.Item(10,3).Value = "Implement 1 cm, 2 cm, and 5 cm length bins in the length data"

$ws.Cells.Item(11,2).Value = "Aggregation treatment"
$ws.Cells.Item(11,3).Value = "Aggregate length and age data before  ('Pre-expansion') or after ('Post-expansion') length and age expansion"

# ============================================================
# Phase-major formatting: apply one property type at a time
# across ALL ranges before moving to the next property type,
# so intermediate/ghost style records are shared instead of
# multiplied per range.
# ============================================================

# --- Phase 1: Font name (all rows 3-11 use Times New Roman) ---
$ws.Range("B3:C11").Font.Name = "Times New Roman"

# --- Phase 2: Font size ---
$ws.Range("B3:C11").Font.Size = 12

# --- Phase 3: Borders ---
$ws.Range("B2:C2").Borders.Item(9).LineStyle = 1
$ws.Range("B3:C3").Borders.Item(8).LineStyle = 1
$ws.Range("B3:C3").Borders.Item(9).LineStyle = 1
$ws.Range("B7:C7").Borders.Item(9).LineStyle = 1
$ws.Range("B8:C8").Borders.Item(9).LineStyle = 1
$ws.Range("B11:C11").Borders.Item(9).LineStyle = 1

# --- Phase 4: Horizontal alignment ---
$ws.Range("B3:C3").HorizontalAlignment = -4108
$ws.Range("B4:B6").HorizontalAlignment = -4131
$ws.Range("B7").HorizontalAlignment = -4131
$ws.Range("B8:C8").HorizontalAlignment = -4108
$ws.Range("B9:B11").HorizontalAlignment = -4131

# --- Phase 5: Vertical alignment ---
$ws.Range("B4:B6").VerticalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108
$ws.Range("B8:C8").VerticalAlignment = -4108
$ws.Range("B9:B11").VerticalAlignment = -4108

# --- Phase 6: Wrap text ---
$ws.Range("C2").WrapText = $true
$ws.Range("C4:C7").WrapText = $true
$ws.Range("C9:C11").WrapText = $true

# --- Merges ---
$ws.Range("B3:C3").Merge()
$ws.Range("B8:C8").Merge()

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 31.5
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 15.75
$ws.Rows.Item(9).RowHeight = 31.5
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 31.5

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 30.3
$ws.Columns.Item(3).ColumnWidth = 80.8

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("B3:C11").Select()
